# Added temp file management & highlight element
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): add 3 new columns H:J ---
$ws.Range("H1").Value = "OptionFromList"
$ws.Range("I1").Value = "FirstOption"
$ws.Range("J1").Value = "SecondOption"
$ws.Range("H1:J1").Copy($ws.Range("H1:J1"))
$ws.Range("A1:G1").Copy()
$ws.Range("H1:J1").PasteSpecial(-4122)

# widths for the new columns
$ws.Range("H1:J1").ColumnWidth = 19.85546875

# --- Row 3: rename existing TestCase1_validateCheckBox -> TestCase3_validateCheckBox ---
$ws.Range("A3").Value = "TestCase3_validateCheckBox"
$ws.Range("H3").Value = "Check Box Demo"

# copy the formatting of the existing data cell (C3, style "1") onto H3/I3/J3/H2/I2/J2
$ws.Range("C3").Copy($ws.Range("H3"))
$ws.Range("H3").Value = "Check Box Demo"
$ws.Range("C2").Copy($ws.Range("H2:J2"))
$ws.Range("C3").Copy($ws.Range("I3:J3"))

# --- Row 4 (new): TestCase4_validateRadioButton ---
$ws.Range("A3:G3").Copy($ws.Range("A4:G4"))
$ws.Range("A4").Value = "TestCase4_validateRadioButton"
$ws.Range("B4").Value = "Validate Radio Buttons"
$ws.Range("C4").Value = "kw_radiobuttons"
$ws.Range("G4").Value = ""
$ws.Range("H3:J3").Copy($ws.Range("H4:J4"))
$ws.Range("H4").Value = "Radio Buttons Demo"
$ws.Range("I4").Value = ""
$ws.Range("J4").Value = ""
$ws.Rows.Item(4).RowHeight = 30

# --- Row 5 (new): TestCase5_validateRadioButton / dropdown list demo ---
$ws.Range("A4:G4").Copy($ws.Range("A5:G5"))
$ws.Range("A5").Value = "TestCase5_validateRadioButton"
$ws.Range("B5").Value = "Validate Drop Down Lists"
$ws.Range("C5").Value = "kw_ddlist"
$ws.Range("H4:J4").Copy($ws.Range("H5:J5"))
$ws.Range("H5").Value = "Select Dropdown List"
$ws.Range("I5").Value = "Tuesday"
$ws.Range("J5").Value = "Florida#New Jersey#Texas#Washington"
$ws.Range("J2").Copy($ws.Range("J5"))
$ws.Range("J5").Value = "Florida#New Jersey#Texas#Washington"
$ws.Rows.Item(5).RowHeight = 45

# --- sheet view / selection ---
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("H9").Select()

# --- page setup (adds pageSetup element with printer settings) ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
